# Saldo.xlsx edit
#
# Changes applied to the "Export" sheet:
#   1) The 004335251 / EDMUNDO row's Saldo changes from -76637.42 to -3000.
#   2) The following three rows (004222784/RAFAEL, 004452912/BRUNO,
#      004361159/HFR) are removed entirely.
#   3) A new row (005232019 / PEDRO / 3000) is inserted right before the
#      005143579 / GABRIEL row.
#   4) A new row (004480134 / JOSE / 39000) is inserted right before the
#      004368468 / AHMAD row.
#
# Operations are applied from the bottom of the sheet upward so that each
# step can use the original (pre-edit) row numbers without having to track
# offsets introduced by earlier insertions/deletions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update the EDMUNDO row's balance and drop the next 3 rows ---
$ws.Cells.Item(215, 3).Value = -3000
$ws.Range("A216:A218").EntireRow.Delete()

# --- Step 2: insert the PEDRO row right before row 26 (005143579/GABRIEL) ---
$ws.Rows.Item(26).Insert()
$ws.Cells.Item(26, 1).NumberFormat = "@"
$ws.Cells.Item(26, 1).Value = "005232019"
$ws.Cells.Item(26, 2).Value = "PEDRO"
$ws.Cells.Item(26, 3).Value = 3000

# --- Step 3: insert the JOSE row right before row 13 (004368468/AHMAD) ---
$ws.Rows.Item(13).Insert()
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "004480134"
$ws.Cells.Item(13, 2).Value = "JOSE"
$ws.Cells.Item(13, 3).Value = 39000
